# ExcelUtil readData test fixture update:
#   - the "Dinnu@247" password used throughout the LoginValidCredentials sheet
#     is replaced with a new password "din1256jlgr"
#   - the previously-selected cell (B11) is replaced by A7
#
# Excel keeps each password-column cell's original mailto hyperlink target
# ("mailto:Dinnu@247") but records the old text via the hyperlink's
# "display" attribute, and collapses the repeated, identical B3:B13
# hyperlinks into a single ranged entry - exactly what happens below.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Start from a clean hyperlink collection so it can be rebuilt in the exact
# shape/order required (this worksheet-level call removes every hyperlink,
# regardless of which Range it's invoked from).
$ws.Hyperlinks.Delete()

# Column A (email) hyperlinks are untouched content-wise - just rebuilt so
# they keep pointing at the same mailto targets as before.
$a2 = $ws.Range("A2")
$null = $ws.Hyperlinks.Add($a2, "mailto:dineshkumar.icon@gmail.com")
$a2.Style = "Hyperlink"

$b2 = $ws.Range("B2")
$null = $ws.Hyperlinks.Add($b2, "mailto:Dinnu@247", "", $null, "Dinnu@247")
$b2.Style = "Hyperlink"

$a3 = $ws.Range("A3")
$null = $ws.Hyperlinks.Add($a3, "mailto:dineshkumar.icon.dk@gmail.com")
$a3.Style = "Hyperlink"

$a4 = $ws.Range("A4")
$null = $ws.Hyperlinks.Add($a4, "mailto:dineshkumar.icon@gmail.com")
$a4.Style = "Hyperlink"

$a5 = $ws.Range("A5")
$null = $ws.Hyperlinks.Add($a5, "mailto:dineshkumar.icon.dk@gmail.com")
$a5.Style = "Hyperlink"

$a6 = $ws.Range("A6")
$null = $ws.Hyperlinks.Add($a6, "mailto:dineshkumar.icon@gmail.com")
$a6.Style = "Hyperlink"

$a7 = $ws.Range("A7")
$null = $ws.Hyperlinks.Add($a7, "mailto:dineshkumar.icon.dk@gmail.com")
$a7.Style = "Hyperlink"

$a8 = $ws.Range("A8")
$null = $ws.Hyperlinks.Add($a8, "mailto:dineshkumar.icon@gmail.com")
$a8.Style = "Hyperlink"

$a9 = $ws.Range("A9")
$null = $ws.Hyperlinks.Add($a9, "mailto:dineshkumar.icon.dk@gmail.com")
$a9.Style = "Hyperlink"

$a10 = $ws.Range("A10")
$null = $ws.Hyperlinks.Add($a10, "mailto:dineshkumar.icon@gmail.com")
$a10.Style = "Hyperlink"

$a12 = $ws.Range("A12")
$null = $ws.Hyperlinks.Add($a12, "mailto:dineshkumar.icon@gmail.com")
$a12.Style = "Hyperlink"

$a11 = $ws.Range("A11")
$null = $ws.Hyperlinks.Add($a11, "mailto:dineshkumar.icon.dk@gmail.com")
$a11.Style = "Hyperlink"

$a13 = $ws.Range("A13")
$null = $ws.Hyperlinks.Add($a13, "mailto:dineshkumar.icon.dk@gmail.com")
$a13.Style = "Hyperlink"

# Column B (password) rows 3-13 all shared the exact same hyperlink, so they
# collapse into a single ranged hyperlink entry, same as Excel would do.
$b3_13 = $ws.Range("B3:B13")
$null = $ws.Hyperlinks.Add($b3_13, "mailto:Dinnu@247", "", $null, "Dinnu@247")
$b3_13.Style = "Hyperlink"

# Now that the hyperlinks (and their historical "display" text) are in
# place, update the actual password text shown in the cells.
$ws.Range("B2:B13").Value = "din1256jlgr"

# Restore the active selection to A7 (was B11 before the edit).
$null = $ws.Range("A7").Select()
